$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.332.66"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.162.08"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.91"
$ws.Range("E5").Value = "  +2.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.63"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.155.24"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("E9").Value = "  +1.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.36"
$ws.Range("E11").Value = "  +1.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("E13").Value = "  +1.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("E14").Value = "  +4.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.687.62"
$ws.Range("E15").Value = "  -0.65%  "

$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.166.24"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.388.96"
$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.07"
$ws.Range("E20").Value = "  +0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.94"
$ws.Range("E21").Value = "  +0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  -1.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.63"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.10"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("E29").Value = "  +3.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.69"
$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.75"
$ws.Range("E31").Value = "  -2.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.05"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("E35").Value = "  -1.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.18"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0727"
$ws.Range("E38").Value = "  +5.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0390"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.12"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  -0.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "390.43"
$ws.Range("E43").Value = "  -5.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.788.26"
$ws.Range("E44").Value = "  -5.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("E46").Value = "  -0.46%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.11"
$ws.Range("E48").Value = "  -2.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.79"
$ws.Range("E49").Value = "  +1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.01"
$ws.Range("E50").Value = "  -1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.111"
$ws.Range("E51").Value = "  +0.56%  "
